$wb = $excel.ActiveWorkbook

# ---- Table 2 sheet ----
$ws1 = $wb.Worksheets.Item("Table 2")

$ws1.Range("A2").Value = "Hospitalized >7 days, acute COVID-19"
$ws1.Range("A3").Value = "Anti-infectives, acute COVID-19"
$ws1.Range("A4").Value = "Anti-platelet, acute COVID-19"
$ws1.Range("A5").Value = "Anti-coagulatives, acute COVID-19"
$ws1.Range("A6").Value = "Immunosuppression, acute COVID-19"

# ---- Table 3 sheet ----
$ws2 = $wb.Worksheets.Item("Table 3")

$ws2.Range("B1").Value = "2-month FUP"
$ws2.Range("C1").Value = "3-month FUP"
$ws2.Range("D1").Value = "6-month FUP"

$ws2.Range("A2").Value = "CT abnormalities" + [char]10 + "at 180-day visit"
$ws2.Range("A3").Value = "CT Severity Score 1-5" + [char]10 + "at 180-day visit"
$ws2.Range("A4").Value = "CT Severity Score >5" + [char]10 + "at 180-day visit"
$ws2.Range("A5").Value = "Symptoms" + [char]10 + "at 180-day visit"
$ws2.Range("A6").Value = "Lung function impairment" + [char]10 + "at 180-day visit"
